# "rearranged music, added songs"
#
# trends sheet ("trends"): row 4 holds the currently-playing song per
# category (col A..L) plus a new "next song" column P (rows 4..11, one
# per upcoming song in the queue). A few songs were swapped out and new
# ones (jeopardy / mexican / dickunddoof / pat / elephant) were queued.
# Also fixes a longstanding typo ("playfullness" -> "playfulness") in
# row 18 / column A.
#
# Shared-string order matters for a byte-faithful round trip, so cells
# that introduce a brand-new string are written in the same order the
# original authoring tool produced them in (typo fix first, then the
# row 4 / column P cells left-to-right, top-to-bottom).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("trends")

# --- typo fix -------------------------------------------------------
$ws.Range("A18").Value = "playfulness"

# --- row 4: swap in the new songs ------------------------------------
$ws.Range("B4").Value = "jeopardy"
$ws.Range("P4").Value = "mexican"
$ws.Range("D4").Value = "dickunddoof"
$ws.Range("F4").Value = "pat"
$ws.Range("I4").Value = "elephant"
$ws.Range("H4").Value = "java"
$ws.Range("J4").Value = "dickunddoof"
$ws.Range("K4").Value = "spanish-flea"
$ws.Range("L4").Value = "jeopardy"

# --- new "up next" column P, rows 5-11 --------------------------------
$ws.Range("P5").Value = "pat"
$ws.Range("P6").Value = "java"
$ws.Range("P7").Value = "elevator"
$ws.Range("P8").Value = "dickunddoof"
$ws.Range("P9").Value = "spanish-flea"
$ws.Range("P10").Value = "elephant"
$ws.Range("P11").Value = "jeopardy"

# --- cursor / selection bookkeeping (matches the saved view state) ----
$ws2 = $wb.Worksheets.Item("do not touch")
$ws2.Activate()
$ws2.Range("A39").Select()

$ws.Activate()
$ws.Range("L7").Select()
